# 9-15-2020 mesocosm_chiachi titrations (initial and 1 week later)
# Appends two new rows (97, 98) of CRM titration accuracy data to Sheet1,
# extends the shared "% off" formula down through the new rows, and
# updates the sheet's active selection to match the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 97: initial titration (CRM opened 2/17/2020 note) ----
# Copy the date-formatted style from the last existing data row (A96) so the
# new date cell keeps the same number format / style index instead of Excel
# inventing a brand-new custom number format.
$ws.Cells.Item(96, 1).Copy() | Out-Null
$ws.Cells.Item(97, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(97, 1).Value = 43885
$ws.Cells.Item(97, 2).Value = 2217.29149784021
$ws.Cells.Item(97, 3).Value = 2214.7603431623102
$ws.Cells.Item(97, 5).Value = 169
$ws.Cells.Item(97, 6).Value = "CRM opened 2/17/2020 (Silbiger bottle for Dudgeon)"

# ---- Row 98: follow-up titration one week later (new CRM bottle) ----
$ws.Cells.Item(96, 1).Copy() | Out-Null
$ws.Cells.Item(98, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(98, 1).Value = 44089
$ws.Cells.Item(98, 2).Value = 2226.5242913638099
$ws.Cells.Item(98, 3).Value = 2224.4699999999998
$ws.Cells.Item(98, 5).Value = 180
$ws.Cells.Item(98, 6).Value = "CRM opened 9/15/2020"

$ws.Application.CutCopyMode = $false

# ---- Extend the "% off" formula down through the new rows ----
# (set per-cell so each keeps a self-contained, recalculable formula)
$ws.Cells.Item(97, 4).Formula = "=100*(B97-C97)/C97"
$ws.Cells.Item(98, 4).Formula = "=100*(B98-C98)/C98"

# ---- Update the window's active cell/selection to follow the new data ----
$ws.Range("D100").Select() | Out-Null
